$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.158.85"
$ws.Range("E2").Value = "  +3.20%  "

$ws.Range("D3").Value = "2.350.05"
$ws.Range("E3").Value = "  +2.33%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "520.62"
$ws.Range("E5").Value = "  +2.67%  "

$ws.Range("E6").Value = "  +3.77%  "

$ws.Range("E7").Value = "  +0.50%  "

$ws.Range("E8").Value = "  +1.29%  "

$ws.Range("D9").Value = "2.348.10"
$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("E10").Value = "  +5.38%  "

$ws.Range("E11").Value = "  -0.80%  "

$ws.Range("E12").Value = "  +3.87%  "

$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").Value = "23.92"
$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").Value = "2.756.59"
$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("D16").Value = "56.986.90"
$ws.Range("E16").Value = "  +3.23%  "

$ws.Range("E17").Value = "  +2.16%  "

$ws.Range("D18").Value = "2.332.36"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("D20").Value = "4.22"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("D21").Value = "323.24"
$ws.Range("E21").Value = "  +4.17%  "

$ws.Range("D22").Value = "6.64"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").Value = "60.89"
$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("E25").Value = "  +8.79%  "

$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("E27").Value = "  +5.46%  "

$ws.Range("D28").Value = "1.31"
$ws.Range("E28").Value = "  +14.07%  "

$ws.Range("E29").Value = "  +5.22%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "168.59"
$ws.Range("E30").Value = "  -2.56%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.71"
$ws.Range("E31").Value = "  +4.63%  "

$ws.Range("E32").Value = "  +0.53%  "

$ws.Range("D33").Value = "18.34"

$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("D37").Value = "0.927"
$ws.Range("E37").Value = "  +0.62%  "

$ws.Range("E38").Value = "  +3.31%  "

$ws.Range("E39").Value = "  +7.40%  "

$ws.Range("D40").Value = "37.93"
$ws.Range("E40").Value = "  +2.96%  "

$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "138.93"
$ws.Range("E42").Value = "  +2.99%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.60"
$ws.Range("E43").Value = "  +4.42%  "

$ws.Range("D44").Value = "5.21"
$ws.Range("E44").Value = "  +5.36%  "

$ws.Range("D45").Value = "277.44"
$ws.Range("E45").Value = "  +5.82%  "

$ws.Range("D46").Value = "0.0933"
$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("E48").Value = "  +1.99%  "

$ws.Range("E49").Value = "  +3.52%  "

$ws.Range("D50").Value = "17.85"
$ws.Range("E50").Value = "  +7.30%  "

$ws.Range("E51").Value = "  +0.41%  "
